# Update the three "HC Collateral ..." headers on the "fact risk" sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("fact risk")

$ws.Range("U1").Value = "HC Collateral Land & Building"
$ws.Range("V1").Value = "HC Collateral Cash, Gold & Other Riskfree Assests"
$ws.Range("W1").Value = "HC Collateral Shares & Other Paper Assests"

# Make "fact risk" the active sheet/tab and move the selection to U1,
# matching the saved view state in the workbook.
$ws.Activate() | Out-Null
$ws.Range("U1").Select() | Out-Null
